# Add data for 2022-08-25: update the "through 08-16" snapshot to "through 08-17"
# (sheet name, label cell, and the August / Total rows of carjacking counts).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "as-of" date.
$ws.Name = "Through 2022-08-17"

# Update the August row (row 9): label + year columns B..I (2015..2022).
$ws.Range("A9").Value = "August (through 08-17)"
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 42
$ws.Range("D9").Value = 43
$ws.Range("E9").Value = 27
$ws.Range("F9").Value = 24
$ws.Range("G9").Value = 111
$ws.Range("H9").Value = 95
$ws.Range("I9").Value = 92

# Update the Total row (row 10): year columns B..I (2015..2022).
$ws.Range("B10").Value = 182
$ws.Range("C10").Value = 344
$ws.Range("D10").Value = 508
$ws.Range("E10").Value = 452
$ws.Range("F10").Value = 328
$ws.Range("G10").Value = 732
$ws.Range("H10").Value = 1005
$ws.Range("I10").Value = 1063

Write-Output ("Sheet name: " + $ws.Name)
Write-Output ("A9: " + $ws.Range("A9").Text)
Write-Output ("Row9: " + $ws.Range("B9").Text + "," + $ws.Range("C9").Text + "," + $ws.Range("D9").Text + "," + $ws.Range("E9").Text + "," + $ws.Range("F9").Text + "," + $ws.Range("G9").Text + "," + $ws.Range("H9").Text + "," + $ws.Range("I9").Text)
Write-Output ("Row10: " + $ws.Range("B10").Text + "," + $ws.Range("C10").Text + "," + $ws.Range("D10").Text + "," + $ws.Range("E10").Text + "," + $ws.Range("F10").Text + "," + $ws.Range("G10").Text + "," + $ws.Range("H10").Text + "," + $ws.Range("I10").Text)
